$wb = $excel.ActiveWorkbook

# --- delta sheet: add column C with constant calibration values (0.1) ---
$wsDelta = $wb.Worksheets.Item("delta")
for ($r = 2; $r -le 20; $r++) {
    $wsDelta.Cells.Item($r, 3).Value = 0.1
}

# --- gamma sheet: add column C with constant calibration values (0.3) ---
$wsGamma = $wb.Worksheets.Item("gamma")
for ($r = 2; $r -le 20; $r++) {
    $wsGamma.Cells.Item($r, 3).Value = 0.3
}

# --- epsilon sheet: add column C with constant calibration values (0.6) ---
$wsEpsilon = $wb.Worksheets.Item("epsilon")
for ($r = 2; $r -le 20; $r++) {
    $wsEpsilon.Cells.Item($r, 3).Value = 0.6
}

# --- alpha sheet: add column C with formula =1/17 (C2 entered alone, then C3:C18 filled together) ---
$wsAlpha = $wb.Worksheets.Item("alpha")
$wsAlpha.Range("C2").Formula = "=1/17"
$wsAlpha.Range("C3:C18").Formula = "=1/17"

# --- update selections on sheets that were visited but not edited (pi) ---
$wsPi = $wb.Worksheets.Item("pi")
[void]$wsPi.Range("C24").Select()

[void]$wsGamma.Range("E20").Select()
[void]$wsEpsilon.Range("C24").Select()
[void]$wsAlpha.Range("C20").Select()

# --- delta ends up the active sheet/tab, with C2:C20 selected ---
[void]$wsDelta.Range("C2:C20").Select()
